$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates -- these are text cells that happen to look like
# numbers, so a leading apostrophe is used to force text entry (mirrors how
# Excel treats quote-prefixed input) and preserve exact formatting
# (trailing zeros, significant digits, etc.)
$ws.Range("D2").Value  = "'261.09"
$ws.Range("D5").Value  = "'0.06096"
$ws.Range("D6").Value  = "'3.513"
$ws.Range("D7").Value  = "'6.719"
$ws.Range("D8").Value  = "'1.361"
$ws.Range("D9").Value  = "'0.7985"
$ws.Range("D10").Value = "'0.1580"
$ws.Range("D11").Value = "'0.08098"
$ws.Range("D12").Value = "'0.03352"
$ws.Range("D13").Value = "'0.03115"
$ws.Range("D15").Value = "'3.893"
$ws.Range("D16").Value = "'0.001687"
$ws.Range("D17").Value = "'0.04830"
$ws.Range("D18").Value = "'0.0006158"
$ws.Range("D19").Value = "'0.006175"
$ws.Range("D20").Value = "'0.001100"
$ws.Range("D21").Value = "'0.003394"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D24").Value = "'2.261"
$ws.Range("D25").Value = "'0.3361"
$ws.Range("D27").Value = "'0.0006163"
$ws.Range("D40").Value = "'0.04596"

# Rows 41-43: the coin list shuffled -- row 41 (BKEXToken) became KickToken,
# row 42 (CEJI) became BKEXToken, row 43 (KickToken) became CEJI, each with
# updated links/prices/labels.
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007091"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1122"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003131"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "'0.01022"
$ws.Range("D46").Value = "'0.00006020"
$ws.Range("D49").Value = "'0.1231"
